# Update cryptos list with refreshed prices / volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.717.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.58%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.095.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.96%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.30%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.46%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.438'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.39'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.64%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.109'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.47%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.380'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.56%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.613.78'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.79%  '

# Row 13
$ws.Range("E13").Value = '  +1.35%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.20%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000166'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.94%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.669.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.54%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.20'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.55%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.090.29'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.10%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.67%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '337.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.64%  '

# Row 22
$ws.Range("E22").Value = '  -0.12%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.508'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.00%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.49%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.173'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.63%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.31%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0919'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.37%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.73%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.30%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.59%  '

# Row 31
$ws.Range("E31").Value = '  +2.67%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.37%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '155.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.20%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.09%  '

# Row 35
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.76%  '

# Row 36
$ws.Range("B36").Value = 'EnergySwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '27.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.56%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.80%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0681'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.60%  '

# Row 39
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.04%  '

# Row 40
$ws.Range("B40").Value = 'RenzoRestakedETH'
$ws.Range("C40").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.125.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.84%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.36%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.31%  '

# Row 43
$ws.Range("E43").Value = '  -1.57%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.286.73'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.12%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.00%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0256'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.59%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.88'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.81%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.955'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.59%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.68%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.737'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.28%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '259.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.27%  '
